{"js": "// Update the 25 two-digit multiplication prompts in the table cells.\n// Each \"NNxNN=\" string is unique within the document, so a targeted\n// search + surgical text replace keeps every other run (fonts, sizes,\n// paragraph marks, empty practice rows, etc.) completely untouched.\nconst replacements = [\n  [\"28\u00d746=\", \"21\u00d757=\"],\n  [\"65\u00d761=\", \"85\u00d791=\"],\n  [\"99\u00d798=\", \"24\u00d753=\"],\n  [\"82\u00d733=\", \"79\u00d717=\"],\n  [\"53\u00d727=\", \"28\u00d793=\"],\n  [\"29\u00d740=\", \"37\u00d786=\"],\n  [\"81\u00d761=\", \"32\u00d719=\"],\n  [\"44\u00d730=\", \"62\u00d774=\"],\n  [\"51\u00d719=\", \"93\u00d793=\"],\n  [\"44\u00d772=\", \"77\u00d797=\"],\n  [\"56\u00d762=\", \"73\u00d718=\"],\n  [\"37\u00d793=\", \"50\u00d785=\"],\n  [\"73\u00d765=\", \"34\u00d779=\"],\n  [\"65\u00d794=\", \"84\u00d778=\"],\n  [\"84\u00d723=\", \"99\u00d742=\"],\n  [\"65\u00d777=\", \"51\u00d727=\"],\n  [\"74\u00d721=\", \"94\u00d786=\"],\n  [\"30\u00d733=\", \"92\u00d735=\"],\n  [\"78\u00d756=\", \"54\u00d790=\"],\n  [\"21\u00d728=\", \"71\u00d796=\"],\n  [\"24\u00d742=\", \"45\u00d737=\"],\n  [\"15\u00d780=\", \"43\u00d764=\"],\n  [\"73\u00d711=\", \"35\u00d774=\"],\n  [\"73\u00d724=\", \"55\u00d729=\"],\n  [\"51\u00d749=\", \"86\u00d797=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find expected text: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 two-digit multiplication prompts in the table cells.\n# Each \"NNxNN=\" string is unique in the document, so Find/Replace on the\n# whole-document range is surgical and leaves every other run (fonts,\n# sizes, paragraph marks, empty practice rows, etc.) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"28\u00d746=\", \"21\u00d757=\"),\n    @(\"65\u00d761=\", \"85\u00d791=\"),\n    @(\"99\u00d798=\", \"24\u00d753=\"),\n    @(\"82\u00d733=\", \"79\u00d717=\"),\n    @(\"53\u00d727=\", \"28\u00d793=\"),\n    @(\"29\u00d740=\", \"37\u00d786=\"),\n    @(\"81\u00d761=\", \"32\u00d719=\"),\n    @(\"44\u00d730=\", \"62\u00d774=\"),\n    @(\"51\u00d719=\", \"93\u00d793=\"),\n    @(\"44\u00d772=\", \"77\u00d797=\"),\n    @(\"56\u00d762=\", \"73\u00d718=\"),\n    @(\"37\u00d793=\", \"50\u00d785=\"),\n    @(\"73\u00d765=\", \"34\u00d779=\"),\n    @(\"65\u00d794=\", \"84\u00d778=\"),\n    @(\"84\u00d723=\", \"99\u00d742=\"),\n    @(\"65\u00d777=\", \"51\u00d727=\"),\n    @(\"74\u00d721=\", \"94\u00d786=\"),\n    @(\"30\u00d733=\", \"92\u00d735=\"),\n    @(\"78\u00d756=\", \"54\u00d790=\"),\n    @(\"21\u00d728=\", \"71\u00d796=\"),\n    @(\"24\u00d742=\", \"45\u00d737=\"),\n    @(\"15\u00d780=\", \"43\u00d764=\"),\n    @(\"73\u00d711=\", \"35\u00d774=\"),\n    @(\"73\u00d724=\", \"55\u00d729=\"),\n    @(\"51\u00d749=\", \"86\u00d797=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #         Format, ReplaceWith, Replace)\n    # Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceAll\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
